# Auto-generated edit script: refresh market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the per-job Leve tables, matching a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 150.25
$ws.Range("I2").Value = 150.25
$ws.Range("K2").Value = 150.25
$ws.Range("M2").Value = -37.25
$ws.Range("H132").Value = 2261.5
$ws.Range("I132").Value = 2261.5
$ws.Range("K132").Value = 6784.5
$ws.Range("M132").Value = -4254.5
$ws.Range("H135").Value = 1641
$ws.Range("I135").Value = 613.1667
$ws.Range("J135").Value = 4724.5
$ws.Range("K135").Value = 5518.5003
$ws.Range("L135").Value = 42520.5
$ws.Range("M135").Value = -2983.5003
$ws.Range("N135").Value = -47590.5
$ws.Range("H137").Value = 1920.75
$ws.Range("I137").Value = 1981.579
$ws.Range("J137").Value = 1792.3334
$ws.Range("K137").Value = 5944.737
$ws.Range("L137").Value = 5377.0002
$ws.Range("M137").Value = -3394.737
$ws.Range("N137").Value = -10477.0002
$ws.Range("H138").Value = 5119.9375
$ws.Range("I138").Value = 2599
$ws.Range("J138").Value = 5701.6924
$ws.Range("K138").Value = 7797
$ws.Range("L138").Value = 17105.0772
$ws.Range("M138").Value = -2657
$ws.Range("N138").Value = -27385.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5033.75
$ws.Range("I2").Value = 6340.6665
$ws.Range("J2").Value = 1113
$ws.Range("K2").Value = 6340.6665
$ws.Range("L2").Value = 1113
$ws.Range("M2").Value = -6227.6665
$ws.Range("N2").Value = -1339
$ws.Range("H32").Value = 5046.4062
$ws.Range("I32").Value = 5046.4062
$ws.Range("K32").Value = 5046.4062
$ws.Range("M32").Value = -4759.4062
$ws.Range("H61").Value = 2383.3333
$ws.Range("I61").Value = 1860
$ws.Range("K61").Value = 1860
$ws.Range("M61").Value = -1648
$ws.Range("H74").Value = 9875.349
$ws.Range("I74").Value = 10254
$ws.Range("K74").Value = 10254
$ws.Range("M74").Value = -9380
$ws.Range("H77").Value = 9875.349
$ws.Range("I77").Value = 10254
$ws.Range("K77").Value = 51270
$ws.Range("M77").Value = -46902
$ws.Range("H116").Value = 5033.75
$ws.Range("I116").Value = 6340.6665
$ws.Range("J116").Value = 1113
$ws.Range("K116").Value = 6340.6665
$ws.Range("L116").Value = 1113
$ws.Range("M116").Value = -4046.6665
$ws.Range("N116").Value = -5701
$ws.Range("H136").Value = 2383.3333
$ws.Range("I136").Value = 1860
$ws.Range("K136").Value = 5580
$ws.Range("M136").Value = -3030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5033.75
$ws.Range("I3").Value = 6340.6665
$ws.Range("J3").Value = 1113
$ws.Range("K3").Value = 6340.6665
$ws.Range("L3").Value = 1113
$ws.Range("M3").Value = -6226.6665
$ws.Range("N3").Value = -1341
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("N10").Value = 0
$ws.Range("H94").Value = 2188.2222
$ws.Range("I94").Value = 2882.5
$ws.Range("J94").Value = 799.6667
$ws.Range("K94").Value = 2882.5
$ws.Range("L94").Value = 799.6667
$ws.Range("M94").Value = -2431.5
$ws.Range("N94").Value = -1701.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 100.25
$ws.Range("H31").Value = 1932.4445
$ws.Range("I31").Value = 1811.625
$ws.Range("J31").Value = 2899
$ws.Range("K31").Value = 1811.625
$ws.Range("L31").Value = 2899
$ws.Range("M31").Value = -1516.625
$ws.Range("N31").Value = -3489
$ws.Range("H34").Value = 1932.4445
$ws.Range("I34").Value = 1811.625
$ws.Range("J34").Value = 2899
$ws.Range("K34").Value = 1811.625
$ws.Range("L34").Value = 2899
$ws.Range("M34").Value = -1609.625
$ws.Range("N34").Value = -3303
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("N42").Value = 0
$ws.Range("H59").Value = 60128.5
$ws.Range("J59").Value = 60128.5
$ws.Range("L59").Value = 60128.5
$ws.Range("N59").Value = -62418.5
$ws.Range("H62").Value = 6006
$ws.Range("J62").Value = 6006
$ws.Range("L62").Value = 6006
$ws.Range("N62").Value = -7254
$ws.Range("H65").Value = 6006
$ws.Range("J65").Value = 6006
$ws.Range("L65").Value = 30030
$ws.Range("N65").Value = -36270
$ws.Range("H132").Value = 3073.4614
$ws.Range("I132").Value = 2878
$ws.Range("K132").Value = 8634
$ws.Range("M132").Value = -6104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H12").Value = 175.88235
$ws.Range("J12").Value = 156.9
$ws.Range("L12").Value = 470.7
$ws.Range("N12").Value = -816.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 443.45456
$ws.Range("I2").Value = 461.2857
$ws.Range("K2").Value = 461.2857
$ws.Range("M2").Value = -348.2857
$ws.Range("H11").Value = 2384653.5
$ws.Range("I11").Value = 2727309
$ws.Range("K11").Value = 2727309
$ws.Range("M11").Value = -2727170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1630.8
$ws.Range("I7").Value = 1813.5
$ws.Range("J7").Value = 900
$ws.Range("K7").Value = 1813.5
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = -1701.5
$ws.Range("N7").Value = -1124
$ws.Range("H46").Value = 2096.2
$ws.Range("I46").Value = 1950
$ws.Range("J46").Value = 2681
$ws.Range("K46").Value = 1950
$ws.Range("L46").Value = 2681
$ws.Range("M46").Value = -1762
$ws.Range("N46").Value = -3057
$ws.Range("H126").Value = 1630.8
$ws.Range("I126").Value = 1813.5
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 5440.5
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -2970.5
$ws.Range("N126").Value = -7640
$ws.Range("H132").Value = 2065.318
$ws.Range("I132").Value = 1841.7333
$ws.Range("K132").Value = 5525.199900000001
$ws.Range("M132").Value = -2995.199900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 803
$ws.Range("I132").Value = 803.75
$ws.Range("K132").Value = 2411.25
$ws.Range("M132").Value = 118.75
$ws.Range("H136").Value = 3530.842
$ws.Range("I136").Value = 2893.6667
$ws.Range("K136").Value = 8681.000100000001
$ws.Range("M136").Value = -6131.000100000001
